$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.00449775112443778
$ws.Cells.Item(2, 3).Value = 0.00224887556221889
$ws.Cells.Item(2, 4).Value = 0.00524737631184408
$ws.Cells.Item(2, 5).Value = 0.00224887556221889
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0.991754122938531
$ws.Cells.Item(2, 8).Value = 0.0217391304347826
$ws.Cells.Item(2, 9).Value = 0.995502248875562
$ws.Cells.Item(2, 10).Value = 0.0194902548725637
$ws.Cells.Item(2, 11).Value = 0.00449775112443778
$ws.Cells.Item(2, 12).Value = 0.0157421289355322
$ws.Cells.Item(2, 13).Value = 0.000749625187406297
$ws.Cells.Item(2, 14).Value = 0.027736131934033
$ws.Cells.Item(2, 15).Value = 0.00149925037481259
$ws.Cells.Item(2, 16).Value = 0.991754122938531
$ws.Cells.Item(2, 17).Value = 0.00524737631184408
$ws.Cells.Item(2, 18).Value = 0.992503748125937
$ws.Cells.Item(2, 19).Value = 0.0434782608695652
$ws.Cells.Item(2, 20).Value = 0.991754122938531
$ws.Cells.Item(2, 21).Value = 0.991004497751124
$ws.Cells.Item(2, 22).Value = 0.00149925037481259
$ws.Cells.Item(2, 23).Value = 0.00299850074962519
$ws.Cells.Item(2, 24).Value = 0.986506746626687
$ws.Cells.Item(3, 2).Value = 0.986506746626687
$ws.Cells.Item(3, 3).Value = 0.992503748125937
$ws.Cells.Item(3, 4).Value = 0.991754122938531
$ws.Cells.Item(3, 5).Value = 0.00224887556221889
$ws.Cells.Item(3, 6).Value = 0.997751124437781
$ws.Cells.Item(3, 7).Value = 0.00149925037481259
$ws.Cells.Item(3, 8).Value = 0.00224887556221889
$ws.Cells.Item(3, 9).Value = 0.00374812593703148
$ws.Cells.Item(3, 10).Value = 0.00224887556221889
$ws.Cells.Item(3, 11).Value = 0.0217391304347826
$ws.Cells.Item(3, 12).Value = 0.979010494752624
$ws.Cells.Item(3, 13).Value = 0.00224887556221889
$ws.Cells.Item(3, 14).Value = 0.00299850074962519
$ws.Cells.Item(3, 15).Value = 0.991754122938531
$ws.Cells.Item(3, 16).Value = 0.00224887556221889
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0.00149925037481259
$ws.Cells.Item(3, 19).Value = 0.00449775112443778
$ws.Cells.Item(3, 20).Value = 0.00149925037481259
$ws.Cells.Item(3, 21).Value = 0.00599700149925037
$ws.Cells.Item(3, 22).Value = 0.00449775112443778
$ws.Cells.Item(3, 23).Value = 0.989505247376312
$ws.Cells.Item(3, 24).Value = 0.00299850074962519
$ws.Cells.Item(4, 2).Value = 0.00374812593703148
$ws.Cells.Item(4, 3).Value = 0.00449775112443778
$ws.Cells.Item(4, 4).Value = 0.00224887556221889
$ws.Cells.Item(4, 5).Value = 0.00149925037481259
$ws.Cells.Item(4, 6).Value = 0.000749625187406297
$ws.Cells.Item(4, 7).Value = 0.00299850074962519
$ws.Cells.Item(4, 8).Value = 0.971514242878561
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0.00449775112443778
$ws.Cells.Item(4, 11).Value = 0.97376311844078
$ws.Cells.Item(4, 12).Value = 0.00224887556221889
$ws.Cells.Item(4, 13).Value = 0.00149925037481259
$ws.Cells.Item(4, 14).Value = 0.968515742128936
$ws.Cells.Item(4, 15).Value = 0.00449775112443778
$ws.Cells.Item(4, 16).Value = 0.00374812593703148
$ws.Cells.Item(4, 17).Value = 0.00299850074962519
$ws.Cells.Item(4, 18).Value = 0.00449775112443778
$ws.Cells.Item(4, 19).Value = 0.00149925037481259
$ws.Cells.Item(4, 20).Value = 0.00224887556221889
$ws.Cells.Item(4, 21).Value = 0.00224887556221889
$ws.Cells.Item(4, 22).Value = 0.988755622188906
$ws.Cells.Item(4, 23).Value = 0.00374812593703148
$ws.Cells.Item(4, 24).Value = 0.00524737631184408
$ws.Cells.Item(5, 2).Value = 0.00524737631184408
$ws.Cells.Item(5, 3).Value = 0.000749625187406297
$ws.Cells.Item(5, 4).Value = 0.000749625187406297
$ws.Cells.Item(5, 5).Value = 0.99400299850075
$ws.Cells.Item(5, 6).Value = 0.00149925037481259
$ws.Cells.Item(5, 7).Value = 0.00374812593703148
$ws.Cells.Item(5, 8).Value = 0.00449775112443778
$ws.Cells.Item(5, 9).Value = 0.000749625187406297
$ws.Cells.Item(5, 10).Value = 0.97376311844078
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0.00224887556221889
$ws.Cells.Item(5, 13).Value = 0.995502248875562
$ws.Cells.Item(5, 14).Value = 0.000749625187406297
$ws.Cells.Item(5, 15).Value = 0.00224887556221889
$ws.Cells.Item(5, 16).Value = 0.00224887556221889
$ws.Cells.Item(5, 17).Value = 0.991754122938531
$ws.Cells.Item(5, 18).Value = 0.00149925037481259
$ws.Cells.Item(5, 19).Value = 0.950524737631184
$ws.Cells.Item(5, 20).Value = 0.00449775112443778
$ws.Cells.Item(5, 21).Value = 0.000749625187406297
$ws.Cells.Item(5, 22).Value = 0.00524737631184408
$ws.Cells.Item(5, 23).Value = 0.00374812593703148
$ws.Cells.Item(5, 24).Value = 0.00524737631184408
